$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Core-Constants")
$ws.Columns.Item(4).Insert()
$ws.Range("D1").Value = "IS Paterne"
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
